# Planning.xlsx update
# - "Revue de littérature" status/progress text updates
# - selection moved to D7
# - columns C and E re-sized (bestfit) to accommodate the longer text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (RECHERCHE line): status of the "REDACTION" column becomes more precise
$ws.Range("C4").Value = "Rédaction finie - 20 pages"

# Row 5 (REVUE DE LITTÉRATURE line)
$ws.Range("E5").Value = "Rédaction finie - 20 pages"

# Row 6 (ENTRETIENS line)
$ws.Range("E6").Value = "Rédaction finie - 10 pages"
$ws.Range("D6").Value = "1/4 rédigé"

# Row 5 continued
$ws.Range("D5").Value = "Entretiens finis + 1/2 rédigé"

# Row 7 (RENDU line) : review milestone date correction 16/3 -> 15/3
$ws.Range("C7").Value = "15/3 : Revue complète"

# Columns C and E grow to fit the new, longer text (keeps the bestFit autosizing behaviour)
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(5).AutoFit()

# Move/leave the active selection on D7, matching the saved cursor position
$ws.Range("D7").Select()
